$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "Buffer", "BufferTime", "ThrottleFirst", "ThrottleLast", "Sample" rows
# as Done (column A) - buffer, bufferTime, bufferCount, bufferTimeOrCount
$ws.Range("A12").Value = "Done"
$ws.Range("A17").Value = "Done"
$ws.Range("A25").Value = "Done"
$ws.Range("A26").Value = "Done"
$ws.Range("A27").Value = "Done"

# Update the active selection/view to reflect the last edited cell
$null = $ws.Range("A12").Select()
